$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 537, shifting the existing rows 537:612 down to 538:613.
$ws.Rows("537:537").Insert()

# Populate the newly inserted row 537 with the same record as the (now shifted)
# row below it, except for the Fecha (D) and Volumen (J) values.
$ws.Range("A537").Value = 5
$ws.Range("B537").Value = "Macroferia Regional de Talca"
$ws.Range("C537").Value = "Maule"
$ws.Range("D537").Value = 45154
$ws.Range("E537").Value = 7
$ws.Range("F537").Value = 100112023
$ws.Range("G537").Value = "Brócoli"
$ws.Range("H537").Value = "Sin especificar"
$ws.Range("I537").Value = "Primera"
$ws.Range("J537").Value = 5000
$ws.Range("K537").Value = 700
$ws.Range("L537").Value = 700
$ws.Range("M537").Value = 700
$ws.Range("N537").Value = "$/unidad"
$ws.Range("O537").Value = "Región del Maule"
$ws.Range("P537").Value = 700
$ws.Range("Q537").Value = 1
$ws.Range("R537").Value = "Hortaliza"
